# Update countries & provincias Spain
# Applies the data refresh described by the diff:
#  - Update the "Datos actualizados..." timestamp string
#  - Swap the Namibia / Principado de Andorra rows (label + data move together)
#  - Swap the Groenlandia / Islas Malvinas rows (label + data move together)
#  - Update the numeric counters (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the timestamp header in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 14:40"

# --- helper to build a 1 x 8 object array for a full data row ---
function New-Row8 {
    param($a,$b,$c,$d,$e,$f,$g,$h)
    $arr = New-Object 'object[,]' 1,8
    $arr[0,0] = $a
    $arr[0,1] = $b
    $arr[0,2] = $c
    $arr[0,3] = $d
    $arr[0,4] = $e
    $arr[0,5] = $f
    $arr[0,6] = $g
    $arr[0,7] = $h
    return $arr
}

# --- 2. Namibia / Principado de Andorra swap (rows 147-148) ---
# Row 147 used to be "Principado de Andorra", row 148 used to be "Namibia".
# After the update row 147 becomes "Namibia" and row 148 becomes
# "Principado de Andorra", each carrying its refreshed counters.
$ws.Range("A147:H147").Value = New-Row8 "Namibia" 861 76 28 832 0 0 1
$ws.Range("A148:H148").Value = New-Row8 "Principado de Andorra" 855 0 803 0 0 0 52

# --- 3. Groenlandia / Islas Malvinas swap (rows 209-210) ---
# Row 209 used to be "Islas Malvinas", row 210 used to be "Groenlandia".
# After the update row 209 becomes "Groenlandia" and row 210 becomes
# "Islas Malvinas" (counters are identical for both, so only the labels move).
$ws.Range("A209:H209").Value = New-Row8 "Groenlandia" 13 0 13 0 0 0 0
$ws.Range("A210:H210").Value = New-Row8 "Islas Malvinas" 13 0 13 0 0 0 0

# --- 4. Refresh numeric counters for the remaining countries ---
# These rows keep their existing country label (column A) and only the
# counters in columns B:H (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) are refreshed.
function Set-Counters {
    param($range, $b,$c,$d,$e,$f,$g,$h)
    $arr = New-Object 'object[,]' 1,7
    $arr[0,0] = $b
    $arr[0,1] = $c
    $arr[0,2] = $d
    $arr[0,3] = $e
    $arr[0,4] = $f
    $arr[0,5] = $g
    $arr[0,6] = $h
    $range.Value = $arr
}

Set-Counters $ws.Range("B4:H4")   3414201 206  1517567 1758839 0 13 137795
Set-Counters $ws.Range("B6:H6")   881846  2380 555992  302617  0 50 23237
Set-Counters $ws.Range("B19:H19") 199998  48   185100  5763    0 1  9135
Set-Counters $ws.Range("B32:H32") 65114   182  55492   9154    0 4  468
Set-Counters $ws.Range("B51:H51") 32941   0    28425   4407    0 1  109
Set-Counters $ws.Range("B68:H68") 13360   363  7852    5446    0 2  62
Set-Counters $ws.Range("B70:H70") 13037   91   12130   297     0 1  610
Set-Counters $ws.Range("B78:H78") 8981    0    8138    590     0 1  253
Set-Counters $ws.Range("B82:H82") 8075    42   3620    4265    0 1  190
Set-Counters $ws.Range("B101:H101") 3775 53   2514    1142    0 0  119
Set-Counters $ws.Range("B110:H110") 2631 14   1981    639     0 0  11
Set-Counters $ws.Range("B119:H119") 1900 4    1871    19      0 0  10
Set-Counters $ws.Range("B169:H169") 245  12   32      211     0 0  2
